$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new user row (row 4) with a base64-encoded password
$ws.Range("A4").Value = "da132224-49b8-4d4d-ba65-e9c25837b300"
$ws.Range("B4").Value = "mustafahere"
$ws.Range("C4").Value = "mustaafhere01@gmail.com"
$ws.Range("D4").Value = "MTIzNDU2"

# Give the new column a custom width, matching the author's save
$ws.Columns.Item(4).ColumnWidth = 13.307291666666666

# Reflect the selection that was active when the workbook was saved
$ws.Range("E7").Select()
